$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 281.6875
$ws.Range("I53").Value = 248.45454
$ws.Range("K53").Value = 248.45454
$ws.Range("M53").Value = 388.54546
$ws.Range("H137").Value = 3708.4167
$ws.Range("I137").Value = 761.8
$ws.Range("K137").Value = 2285.4
$ws.Range("M137").Value = 264.6000000000004
$ws.Range("H138").Value = 4312.8
$ws.Range("J138").Value = 4510.567
$ws.Range("L138").Value = 13531.701
$ws.Range("N138").Value = -23811.701

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 844.125
$ws.Range("I2").Value = 828.375
$ws.Range("K2").Value = 828.375
$ws.Range("M2").Value = -715.375
$ws.Range("H61").Value = 2757.3
$ws.Range("I61").Value = 2619.4443
$ws.Range("K61").Value = 2619.4443
$ws.Range("M61").Value = -2407.4443
$ws.Range("H110").Value = 2179.4666
$ws.Range("I110").Value = 2052.1428
$ws.Range("K110").Value = 2052.1428
$ws.Range("M110").Value = -7.142800000000079
$ws.Range("H116").Value = 844.125
$ws.Range("I116").Value = 828.375
$ws.Range("K116").Value = 828.375
$ws.Range("M116").Value = 1465.625
$ws.Range("H132").Value = 4027.4285
$ws.Range("I132").Value = 2932
$ws.Range("K132").Value = 8796
$ws.Range("M132").Value = -6266
$ws.Range("H136").Value = 2757.3
$ws.Range("I136").Value = 2619.4443
$ws.Range("K136").Value = 7858.3329
$ws.Range("M136").Value = -5308.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 844.125
$ws.Range("I3").Value = 828.375
$ws.Range("K3").Value = 828.375
$ws.Range("M3").Value = -714.375
$ws.Range("H86").Value = 1212.909
$ws.Range("I86").Value = 1212.909
$ws.Range("K86").Value = 1212.909
$ws.Range("M86").Value = -89.90900000000011
$ws.Range("H89").Value = 1212.909
$ws.Range("I89").Value = 1212.909
$ws.Range("K89").Value = 6064.545
$ws.Range("M89").Value = -448.5450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6829.1904
$ws.Range("I31").Value = 2446.6365
$ws.Range("J31").Value = 11650
$ws.Range("K31").Value = 2446.6365
$ws.Range("L31").Value = 11650
$ws.Range("M31").Value = -2151.6365
$ws.Range("N31").Value = -12240
$ws.Range("H34").Value = 6829.1904
$ws.Range("I34").Value = 2446.6365
$ws.Range("J34").Value = 11650
$ws.Range("K34").Value = 2446.6365
$ws.Range("L34").Value = 11650
$ws.Range("M34").Value = -2244.6365
$ws.Range("N34").Value = -12054
$ws.Range("H58").Value = 2206.1
$ws.Range("I58").Value = 2206.1
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2206.1
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2003.1
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 14002
$ws.Range("I86").Value = 14002
$ws.Range("K86").Value = 14002
$ws.Range("M86").Value = -12879
$ws.Range("H89").Value = 14002
$ws.Range("I89").Value = 14002
$ws.Range("K89").Value = 70010
$ws.Range("M89").Value = -64394
$ws.Range("H136").Value = 2206.1
$ws.Range("I136").Value = 2206.1
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6618.299999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4068.299999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 53.555557
$ws.Range("I2").Value = 48.833332
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 292.999992
$ws.Range("L2").Value = 378
$ws.Range("M2").Value = -179.999992
$ws.Range("N2").Value = -604
$ws.Range("H5").Value = 807.7143
$ws.Range("J5").Value = 887.375
$ws.Range("L5").Value = 2662.125
$ws.Range("N5").Value = -2886.125
$ws.Range("H21").Value = 495
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 495
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 1485
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -1831
$ws.Range("H34").Value = 2422.5
$ws.Range("J34").Value = 2945
$ws.Range("L34").Value = 8835
$ws.Range("N34").Value = -9003
$ws.Range("H39").Value = 7000
$ws.Range("J39").Value = 7000
$ws.Range("L39").Value = 21000
$ws.Range("N39").Value = -21588
$ws.Range("H59").Value = 866.6667
$ws.Range("I59").Value = 866.6667
$ws.Range("K59").Value = 2600.0001
$ws.Range("M59").Value = -2060.0001
$ws.Range("H80").Value = 8995.75
$ws.Range("J80").Value = 8993
$ws.Range("L80").Value = 26979
$ws.Range("N80").Value = -28851
$ws.Range("H83").Value = 8995.75
$ws.Range("J83").Value = 8993
$ws.Range("L83").Value = 80937
$ws.Range("N83").Value = -90297
$ws.Range("H109").Value = 1244.8334
$ws.Range("I109").Value = 784.5
$ws.Range("J109").Value = 2165.5
$ws.Range("K109").Value = 2353.5
$ws.Range("L109").Value = 6496.5
$ws.Range("M109").Value = -1313.5
$ws.Range("N109").Value = -8576.5
$ws.Range("H122").Value = 3718.5557
$ws.Range("J122").Value = 3892.762
$ws.Range("L122").Value = 35034.858
$ws.Range("N122").Value = -39934.858
$ws.Range("H131").Value = 1418.6
$ws.Range("J131").Value = 1499.8334
$ws.Range("L131").Value = 4499.5002
$ws.Range("N131").Value = -14579.5002
$ws.Range("H132").Value = 4494.4443
$ws.Range("I132").Value = 2663
$ws.Range("K132").Value = 23967
$ws.Range("M132").Value = -21437
$ws.Range("H135").Value = 807.7143
$ws.Range("J135").Value = 887.375
$ws.Range("L135").Value = 7986.375
$ws.Range("N135").Value = -13056.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3029
$ws.Range("I132").Value = 2604.762
$ws.Range("K132").Value = 7814.286
$ws.Range("M132").Value = -5284.286
$ws.Range("H140").Value = 125000
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4417
$ws.Range("I61").Value = 4403.625
$ws.Range("K61").Value = 4403.625
$ws.Range("M61").Value = -4201.625
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H113").Value = 4417
$ws.Range("I113").Value = 4403.625
$ws.Range("K113").Value = 4403.625
$ws.Range("M113").Value = -2233.625
$ws.Range("H122").Value = 2466.3333
$ws.Range("I122").Value = 2559.6
$ws.Range("K122").Value = 7678.799999999999
$ws.Range("M122").Value = -5228.799999999999
$ws.Range("H132").Value = 4381.0557
$ws.Range("I132").Value = 3948.875
$ws.Range("J132").Value = 4726.8
$ws.Range("K132").Value = 11846.625
$ws.Range("L132").Value = 14180.4
$ws.Range("M132").Value = -9316.625
$ws.Range("N132").Value = -19240.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 1139.8
$ws.Range("I81").Value = 1139.8
$ws.Range("K81").Value = 2279.6
$ws.Range("M81").Value = -1218.6
$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766
$ws.Range("H84").Value = 1139.8
$ws.Range("I84").Value = 1139.8
$ws.Range("K84").Value = 11398
$ws.Range("M84").Value = -6094
$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652
